$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Schistosoma mansoni ---
$ws.Range("G2").Value = 9

# --- Row 3: Brugia malayi ---
$ws.Range("D3").NumberFormat = "d-mmm"
$ws.Range("D3").Value = 43470
$ws.Range("G3").Value = 2

# --- Row 5: Schistosoma haematobium ---
$ws.Range("G5").Value = 11

# --- Row 6: Wuchereria bancrofti ---
$ws.Range("G6").Value = 18

# --- Row 8: Onchocerca volvulus ---
$ws.Range("G8").Value = 6

# --- Row 11: Loa loa ---
$ws.Range("D11").NumberFormat = "d-mmm"
$ws.Range("D11").Value = 43462
$ws.Range("G11").Value = 9

# --- Row 15: Echinococcus granulosus ---
$ws.Range("D15").NumberFormat = "d-mmm"
$ws.Range("D15").Value = 43456
$ws.Range("G15").Value = 1

# --- Row 16: Taenia solium ---
$ws.Range("G16").Value = 8

# --- Row 25: Brugia timori ---
$ws.Range("D25").NumberFormat = "d-mmm"
$ws.Range("D25").Value = 43470
$ws.Range("G25").Value = 2

# --- Row 26: Mansonella perstans ---
$ws.Range("E26").NumberFormat = "d-mmm"
$ws.Range("E26").Value = 43470
$ws.Range("F26").Value = 104
$ws.Range("G26").Value = 5

# --- Row 27: Trichostrongylus colubriformis ---
$ws.Range("D27").NumberFormat = "d-mmm"
$ws.Range("D27").Value = 43462
$ws.Range("E27").NumberFormat = "d-mmm"
$ws.Range("E27").Value = 43462
$ws.Range("G27").Value = 0

# --- Row 28: Taenia saginata ---
$ws.Range("D28").NumberFormat = "d-mmm"
$ws.Range("D28").Value = 43456
$ws.Range("E28").NumberFormat = "d-mmm"
$ws.Range("E28").Value = 43461
$ws.Range("G28").Value = 4

# --- Row 29: Angiostrongylus cantonensis ---
$ws.Range("E29").NumberFormat = "d-mmm"
$ws.Range("E29").Value = 43455
$ws.Range("F29").Value = 90
$ws.Range("G29").Value = 3

# --- Row 30: Enterobius vermicularis ---
$ws.Range("D30").NumberFormat = "d-mmm"
$ws.Range("D30").Value = 43454
$ws.Range("E30").NumberFormat = "d-mmm"
$ws.Range("E30").Value = 43454
$ws.Range("F30").Value = 68
$ws.Range("G30").Value = 0

# --- Row 31: Taenia crassiceps ---
$ws.Range("D31").NumberFormat = "d-mmm"
$ws.Range("D31").Value = 43454
$ws.Range("E31").NumberFormat = "d-mmm"
$ws.Range("E31").Value = 43454
$ws.Range("F31").Value = 63
$ws.Range("G31").Value = 0

# --- Row 32: Paragonimus westermani ---
$ws.Range("D32").NumberFormat = "d-mmm"
$ws.Range("D32").Value = 43454
$ws.Range("E32").NumberFormat = "d-mmm"
$ws.Range("E32").Value = 43454
$ws.Range("F32").Value = 60
$ws.Range("G32").Value = 1

# --- Row 33: Hymenolepis nana ---
$ws.Range("F33").Value = 60

# --- Row 95: Paragonimus skrjabini ---
$ws.Range("G95").Value = 1

# --- Update selection to A25 ---
$ws.Range("A25").Select()
